# Fix the merge-field placeholder "{{ubicación}}" -> "{{ubicacion}}"
# (drop the accent on the "o" so it matches the certificate-generation
# code's variable name), reproducing the run split seen in the target
# revision: {{ubicaci | o | n}} as three separate <w:r> runs that keep
# the original run formatting (Arial Narrow, bold).

$d = $word.ActiveDocument

# Locate the literal ASCII prefix "{{ubicaci" (unique in the document) so
# we never have to depend on Find() matching the accented character.
$prefix = $d.Content
$foundPrefix = $prefix.Find.Execute("{{ubicaci", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($foundPrefix) {
    $prefixStart = $prefix.Start
    $prefixEnd = $prefix.End

    # The accented character ("ó") immediately follows the prefix.
    $accentStart = $prefixEnd
    $accentEnd = $accentStart + 1

    # Replace "ó" with a plain "o" (keeps the field the same length).
    $accentRange = $d.Range($accentStart, $accentEnd)
    $accentRange.Text = "o"

    # Force the edited middle character into its own run (distinct from
    # the unedited prefix/suffix runs) by toggling Bold off/on, matching
    # the three-run structure produced by the original edit.
    $midRange = $d.Range($accentStart, $accentEnd)
    $midRange.Bold = 0
    $midRange.Bold = 1
}
